$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Delete the trailing empty paragraph (the very last paragraph, which has
#    ind left=426 and no text) right before the sectPr.
# ---------------------------------------------------------------------------
$lastIndex = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($lastIndex)
if ($lastPara.Range.Text -eq [char]13) {
    $lastPara.Range.Delete()
}

# ---------------------------------------------------------------------------
# 2) Remove the "_GoBack" bookmark (left over editing artifact).
# ---------------------------------------------------------------------------
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

# ---------------------------------------------------------------------------
# 3) Turn the manual line breaks in the title paragraph into " | " (and a
#    plain " " before the last segment) separators instead of <w:br/>.
# ---------------------------------------------------------------------------
$titlePara = $d.Paragraphs.Item(1)

$sep = " | "
$r = $d.Range($titlePara.Range.Start, $titlePara.Range.End)
$null = $r.Find.Execute([char]11, $false, $false, $false, $false, $false, $true, 0, $false, $sep, 1)

$r = $d.Range($titlePara.Range.Start, $titlePara.Range.End)
$null = $r.Find.Execute([char]11, $false, $false, $false, $false, $false, $true, 0, $false, $sep, 1)

$r = $d.Range($titlePara.Range.Start, $titlePara.Range.End)
$null = $r.Find.Execute([char]11, $false, $false, $false, $false, $false, $true, 0, $false, $sep, 1)

$r = $d.Range($titlePara.Range.Start, $titlePara.Range.End)
$null = $r.Find.Execute([char]11, $false, $false, $false, $false, $false, $true, 0, $false, " ", 1)

# ---------------------------------------------------------------------------
# 4) Delete the old "Author: JJ van Zon" / "Location: ..." block, which sat
#    between the title paragraph and the "Super-Project" heading, along
#    with the two blank paragraphs surrounding it. This causes the existing
#    "Super-Project" heading (further down) to become the paragraph right
#    after the title, matching the target order.
# ---------------------------------------------------------------------------
for ($i = 5; $i -ge 2; $i--) {
    $p = $d.Paragraphs.Item($i)
    $p.Range.Delete()
}

# ---------------------------------------------------------------------------
# 5) Bold the numbers "3" (in "3 weeks") and "52 ¼" (in "52 ¼ hours of
#    work"). Locate the paragraphs by their text content (more robust than
#    a document-wide Find, whose Range can be confusing to re-scope).
# ---------------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text

    if ($t -like "3 weeks*") {
        $numRange = $d.Range($p.Range.Start, $p.Range.Start + 1)
        $numRange.Font.Bold = 1
        $numRange.Font.BoldBi = 1
    }
    elseif ($t -like "*hours of work*") {
        $numRange2 = $d.Range($p.Range.Start, $p.Range.Start + 4)
        $numRange2.Font.Bold = 1
        $numRange2.Font.BoldBi = 1
    }
}

Write-Host "Final paragraph count:" $d.Paragraphs.Count
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    Write-Host $i ":" $p.Range.Text
}
